$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# A new tracking entry was logged on top of the September list (row 30),
# pushing all the existing September/August/row-58 data down by one row.
$ws.Rows(30).Insert()

$ws.Range("R30").Value = "bal axis"
$ws.Range("S30").Value = "2024-09-05 09:06:25"
